# Updates cryptos list figures (price + volume/1h %) per the
# "Updated cryptos list ... with GitHub Actions" commit.
# Rows 19/20 (TRON/Uniswap) and 32/33 (ImmutableX/WrappedeETH) are
# re-ordered in place by overwriting each column in the pair.
#
# Numeric-looking price strings (column D) are written with a
# temporary Text number format so Excel keeps them as strings
# instead of auto-converting to numbers; the format/style is
# reset to Normal right after so no stray formatting is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.776.47'
$ws.Range('E2').Value = '  -0.34%  '

# Row 3
$ws.Range('D3').Value = '3.840.68'
$ws.Range('E3').Value = '  +2.34%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '602.41'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.05%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '162.66'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.83%  '

# Row 7
$ws.Range('D7').Value = '3.833.15'
$ws.Range('E7').Value = '  +2.21%  '

# Row 8
$ws.Range('E8').Value = '  +0.09%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.530'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.77%  '

# Row 10
$ws.Range('E10').Value = '  -0.82%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.29'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.50%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.459'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.23%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '36.82'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.14%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000243'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.07%  '

# Row 15
$ws.Range('D15').Value = '4.480.17'
$ws.Range('E15').Value = '  +2.25%  '

# Row 16
$ws.Range('D16').Value = '3.839.86'
$ws.Range('E16').Value = '  +2.26%  '

# Row 17
$ws.Range('D17').Value = '68.994.74'
$ws.Range('E17').Value = '  -0.08%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.55'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.87%  '

# Row 19
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.46'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +5.15%  '

# Row 20
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.113'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.14%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.11'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.77%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '484.21'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.68%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.717'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.10%  '

# Row 24
$ws.Range('E24').Value = '  +3.48%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.04'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.93%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.25'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.04%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.10'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.84%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.99'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.21%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.16%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.96'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.35%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.91'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.26%  '

# Row 32
$ws.Range('B32').Value = 'WrappedeETH'
$ws.Range('C32').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D32').Value = '3.997.48'
$ws.Range('E32').Value = '  +2.57%  '

# Row 33
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.38'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.06%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '32.19'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.81%  '

# Row 35
$ws.Range('D35').Value = '3.788.67'
$ws.Range('E35').Value = '  +2.76%  '

# Row 36
$ws.Range('E36').Value = '  -1.79%  '

# Row 37
$ws.Range('E37').Value = '  +1.32%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.141'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.16%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.86'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.42%  '

# Row 40
$ws.Range('E40').Value = '  -0.03%  '

# Row 41
$ws.Range('E41').Value = '  -1.92%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '438.49'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.02%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.97'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.82%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '48.48'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.34%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.97'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.53%  '

# Row 46
$ws.Range('E46').Value = '  -0.01%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.38'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.28%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '27.42'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +17.46%  '

# Row 49
$ws.Range('D49').Value = '2.842.25'
$ws.Range('E49').Value = '  +1.64%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '142.60'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.87%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0356'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.13%  '
